$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.545.83'
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").Value = '1.568.67'
$ws.Range("E3").Value = '  -1.57%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.70'
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("E6").Value = '  -0.56%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.18'
$ws.Range("E8").Value = '  +5.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.10'
$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("E10").Value = '  -1.73%  '

$ws.Range("E11").Value = '  -1.75%  '

$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").Value = '1.568.81'
$ws.Range("E14").Value = '  -1.59%  '

$ws.Range("E15").Value = '  -2.07%  '

$ws.Range("D16").Value = '28.520.64'
$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("E17").Value = '  -2.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.20'
$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.32'
$ws.Range("E19").Value = '  +1.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").Value = '  -1.83%  '

$ws.Range("E21").Value = '  -2.69%  '

$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("E23").Value = '  -6.04%  '

$ws.Range("E24").Value = '  -2.44%  '

$ws.Range("E25").Value = '  +8.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.20'
$ws.Range("E26").Value = '  -0.30%  '

$ws.Range("E27").Value = '  -1.20%  '

$ws.Range("E28").Value = '  -2.67%  '

$ws.Range("E29").Value = '  -3.39%  '

$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("E31").Value = '  +1.97%  '

$ws.Range("E32").Value = '  -3.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.20'
$ws.Range("E33").Value = '  -1.16%  '

$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("D35").Value = '1.392.48'
$ws.Range("E35").Value = '  -0.40%  '

$ws.Range("E36").Value = '  +0.71%  '

$ws.Range("E37").Value = '  -3.68%  '

$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("E39").Value = '  +2.80%  '

$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.523'
$ws.Range("E41").Value = '  -3.25%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("E44").Value = '  -3.19%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.48'
$ws.Range("E45").Value = '  -4.14%  '

$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0462'
$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '62.80'
$ws.Range("E48").Value = '  -2.46%  '

$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.35'
$ws.Range("E50").Value = '  -1.28%  '

$ws.Range("E51").Value = '  -0.87%  '
